$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1633.3684
$ws.Cells.Item(18, 9).Value = 1663
$ws.Cells.Item(18, 11).Value = 1663
$ws.Cells.Item(18, 13).Value = -1379
$ws.Cells.Item(20, 8).Value = 1634.5
$ws.Cells.Item(20, 9).Value = 1634.5
$ws.Cells.Item(20, 11).Value = 1634.5
$ws.Cells.Item(20, 13).Value = -1404.5
$ws.Cells.Item(35, 8).Value = 1634.5
$ws.Cells.Item(35, 9).Value = 1634.5
$ws.Cells.Item(35, 11).Value = 1634.5
$ws.Cells.Item(35, 13).Value = -1255.5
$ws.Cells.Item(92, 8).Value = 463.83334
$ws.Cells.Item(92, 9).Value = 414.63635
$ws.Cells.Item(92, 11).Value = 414.63635
$ws.Cells.Item(92, 13).Value = 833.36365
$ws.Cells.Item(96, 8).Value = 310.6154
$ws.Cells.Item(96, 9).Value = 353.9091
$ws.Cells.Item(96, 10).Value = 72.5
$ws.Cells.Item(96, 11).Value = 1061.7273
$ws.Cells.Item(96, 12).Value = 217.5
$ws.Cells.Item(96, 13).Value = 311.2727
$ws.Cells.Item(96, 14).Value = -2963.5
$ws.Cells.Item(100, 8).Value = 1524.5454
$ws.Cells.Item(100, 9).Value = 1271.25
$ws.Cells.Item(100, 10).Value = 2200
$ws.Cells.Item(100, 11).Value = 1271.25
$ws.Cells.Item(100, 12).Value = 2200
$ws.Cells.Item(100, 13).Value = -730.25
$ws.Cells.Item(100, 14).Value = -3282
$ws.Cells.Item(137, 8).Value = 1869.375
$ws.Cells.Item(137, 9).Value = 1650.8422
$ws.Cells.Item(137, 11).Value = 4952.5266
$ws.Cells.Item(137, 13).Value = -2402.5266

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5500.933
$ws.Cells.Item(32, 9).Value = 5500.933
$ws.Cells.Item(32, 11).Value = 5500.933
$ws.Cells.Item(32, 13).Value = -5213.933
$ws.Cells.Item(61, 8).Value = 31259922
$ws.Cells.Item(61, 9).Value = 50008300
$ws.Cells.Item(61, 11).Value = 50008300
$ws.Cells.Item(61, 13).Value = -50008088
$ws.Cells.Item(74, 8).Value = 2577.8462
$ws.Cells.Item(74, 10).Value = 3819.3333
$ws.Cells.Item(74, 12).Value = 3819.3333
$ws.Cells.Item(74, 14).Value = -5567.3333
$ws.Cells.Item(77, 8).Value = 2577.8462
$ws.Cells.Item(77, 10).Value = 3819.3333
$ws.Cells.Item(77, 12).Value = 19096.6665
$ws.Cells.Item(77, 14).Value = -27832.6665
$ws.Cells.Item(132, 8).Value = 3038.0195
$ws.Cells.Item(132, 9).Value = 2571.7292
$ws.Cells.Item(132, 10).Value = 10498.667
$ws.Cells.Item(132, 11).Value = 7715.187600000001
$ws.Cells.Item(132, 12).Value = 31496.001
$ws.Cells.Item(132, 13).Value = -5185.187600000001
$ws.Cells.Item(132, 14).Value = -36556.001
$ws.Cells.Item(136, 8).Value = 31259922
$ws.Cells.Item(136, 9).Value = 50008300
$ws.Cells.Item(136, 11).Value = 150024900
$ws.Cells.Item(136, 13).Value = -150022350

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2492.5
$ws.Cells.Item(86, 9).Value = 2466
$ws.Cells.Item(86, 10).Value = 2519
$ws.Cells.Item(86, 11).Value = 2466
$ws.Cells.Item(86, 12).Value = 2519
$ws.Cells.Item(86, 13).Value = -1343
$ws.Cells.Item(86, 14).Value = -4765
$ws.Cells.Item(89, 8).Value = 2492.5
$ws.Cells.Item(89, 9).Value = 2466
$ws.Cells.Item(89, 10).Value = 2519
$ws.Cells.Item(89, 11).Value = 12330
$ws.Cells.Item(89, 12).Value = 12595
$ws.Cells.Item(89, 13).Value = -6714
$ws.Cells.Item(89, 14).Value = -23827
$ws.Cells.Item(94, 8).Value = 2827.4375
$ws.Cells.Item(94, 9).Value = 2520.0833
$ws.Cells.Item(94, 11).Value = 2520.0833
$ws.Cells.Item(94, 13).Value = -2069.0833
$ws.Cells.Item(134, 8).Value = 2293.4
$ws.Cells.Item(134, 9).Value = 2333.4524
$ws.Cells.Item(134, 11).Value = 7000.3572
$ws.Cells.Item(134, 13).Value = -4465.3572

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4186.5483
$ws.Cells.Item(31, 9).Value = 3375.5881
$ws.Cells.Item(31, 11).Value = 3375.5881
$ws.Cells.Item(31, 13).Value = -3080.5881
$ws.Cells.Item(34, 8).Value = 4186.5483
$ws.Cells.Item(34, 9).Value = 3375.5881
$ws.Cells.Item(34, 11).Value = 3375.5881
$ws.Cells.Item(34, 13).Value = -3173.5881
$ws.Cells.Item(74, 8).Value = 34494.285
$ws.Cells.Item(74, 10).Value = 34494.285
$ws.Cells.Item(74, 12).Value = 34494.285
$ws.Cells.Item(74, 14).Value = -36242.285
$ws.Cells.Item(77, 8).Value = 34494.285
$ws.Cells.Item(77, 10).Value = 34494.285
$ws.Cells.Item(77, 12).Value = 103482.855
$ws.Cells.Item(77, 14).Value = -112218.855
$ws.Cells.Item(86, 8).Value = 7531
$ws.Cells.Item(86, 9).Value = 6546.6665
$ws.Cells.Item(86, 10).Value = 9499.666999999999
$ws.Cells.Item(86, 11).Value = 6546.6665
$ws.Cells.Item(86, 12).Value = 9499.666999999999
$ws.Cells.Item(86, 13).Value = -5423.6665
$ws.Cells.Item(86, 14).Value = -11745.667
$ws.Cells.Item(89, 8).Value = 7531
$ws.Cells.Item(89, 9).Value = 6546.6665
$ws.Cells.Item(89, 10).Value = 9499.666999999999
$ws.Cells.Item(89, 11).Value = 32733.3325
$ws.Cells.Item(89, 12).Value = 47498.335
$ws.Cells.Item(89, 13).Value = -27117.3325
$ws.Cells.Item(89, 14).Value = -58730.335
$ws.Cells.Item(105, 8).Value = 1901.4166
$ws.Cells.Item(105, 9).Value = 2071.7
$ws.Cells.Item(105, 10).Value = 1050
$ws.Cells.Item(105, 11).Value = 2071.7
$ws.Cells.Item(105, 12).Value = 1050
$ws.Cells.Item(105, 13).Value = -324.6999999999998
$ws.Cells.Item(105, 14).Value = -4544
$ws.Cells.Item(112, 8).Value = 75973.5
$ws.Cells.Item(112, 10).Value = 75973.5
$ws.Cells.Item(112, 12).Value = 75973.5
$ws.Cells.Item(112, 14).Value = -78927.5
$ws.Cells.Item(132, 8).Value = 2642.05
$ws.Cells.Item(132, 9).Value = 1793.1333
$ws.Cells.Item(132, 11).Value = 5379.3999
$ws.Cells.Item(132, 13).Value = -2849.3999
$ws.Cells.Item(134, 8).Value = 9225
$ws.Cells.Item(134, 9).Value = 7300
$ws.Cells.Item(134, 11).Value = 21900
$ws.Cells.Item(134, 13).Value = -19365

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 4218
$ws.Cells.Item(122, 9).Value = 7255
$ws.Cells.Item(122, 11).Value = 65295
$ws.Cells.Item(122, 13).Value = -62845

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, 8).Value = 38571.43
$ws.Cells.Item(98, 10).Value = 38571.43
$ws.Cells.Item(98, 12).Value = 38571.43
$ws.Cells.Item(98, 14).Value = -44561.43
$ws.Cells.Item(102, 8).Value = 5308.15
$ws.Cells.Item(102, 9).Value = 4798.0527
$ws.Cells.Item(102, 10).Value = 15000
$ws.Cells.Item(102, 11).Value = 4798.0527
$ws.Cells.Item(102, 12).Value = 15000
$ws.Cells.Item(102, 13).Value = -3176.0527
$ws.Cells.Item(102, 14).Value = -18244
$ws.Cells.Item(132, 8).Value = 3195.9795
$ws.Cells.Item(132, 9).Value = 2747.973
$ws.Cells.Item(132, 10).Value = 4577.3335
$ws.Cells.Item(132, 11).Value = 8243.919
$ws.Cells.Item(132, 12).Value = 13732.0005
$ws.Cells.Item(132, 13).Value = -5713.919
$ws.Cells.Item(132, 14).Value = -18792.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2166.6667
$ws.Cells.Item(27, 8).Value = 2166.6667
$ws.Cells.Item(55, 8).Value = 817.8421
$ws.Cells.Item(55, 10).Value = 630.8889
$ws.Cells.Item(55, 12).Value = 630.8889
$ws.Cells.Item(55, 14).Value = -976.8889
$ws.Cells.Item(122, 8).Value = 3628.4
$ws.Cells.Item(122, 9).Value = 2697.6667
$ws.Cells.Item(122, 11).Value = 8093.000100000001
$ws.Cells.Item(122, 13).Value = -5643.000100000001
$ws.Cells.Item(132, 8).Value = 7386.3335
$ws.Cells.Item(132, 9).Value = 8031.8696
$ws.Cells.Item(132, 11).Value = 24095.6088
$ws.Cells.Item(132, 13).Value = -21565.6088
$ws.Cells.Item(134, 8).Value = 94082
$ws.Cells.Item(134, 9).Value = 93900
$ws.Cells.Item(134, 10).Value = 94142.664
$ws.Cells.Item(134, 11).Value = 93900
$ws.Cells.Item(134, 12).Value = 94142.664
$ws.Cells.Item(134, 13).Value = -88830
$ws.Cells.Item(134, 14).Value = -104282.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 4881
$ws.Cells.Item(4, 9).Value = 3905
$ws.Cells.Item(4, 11).Value = 3905
$ws.Cells.Item(4, 13).Value = -3792
$ws.Cells.Item(12, 8).Value = 10999.667
$ws.Cells.Item(12, 9).Value = 11499.5
$ws.Cells.Item(12, 11).Value = 11499.5
$ws.Cells.Item(12, 13).Value = -11357.5
$ws.Cells.Item(14, 8).Value = 15133.896
$ws.Cells.Item(14, 9).Value = 13520.211
$ws.Cells.Item(14, 10).Value = 18199.9
$ws.Cells.Item(14, 11).Value = 13520.211
$ws.Cells.Item(14, 12).Value = 18199.9
$ws.Cells.Item(14, 13).Value = -13352.211
$ws.Cells.Item(14, 14).Value = -18535.9
$ws.Cells.Item(45, 8).Value = 35498.75
$ws.Cells.Item(45, 9).Value = 35665
$ws.Cells.Item(45, 10).Value = 35000
$ws.Cells.Item(45, 11).Value = 35665
$ws.Cells.Item(45, 12).Value = 35000
$ws.Cells.Item(45, 13).Value = -35174
$ws.Cells.Item(45, 14).Value = -35982
$ws.Cells.Item(126, 8).Value = 2192.6667
$ws.Cells.Item(126, 9).Value = 2291.25
$ws.Cells.Item(126, 10).Value = 1995.5
$ws.Cells.Item(126, 11).Value = 6873.75
$ws.Cells.Item(126, 12).Value = 5986.5
$ws.Cells.Item(126, 13).Value = -4403.75
$ws.Cells.Item(126, 14).Value = -10926.5
$ws.Cells.Item(132, 8).Value = 3092.7036
$ws.Cells.Item(132, 9).Value = 2747.6191
$ws.Cells.Item(132, 11).Value = 8242.8573
$ws.Cells.Item(132, 13).Value = -5712.8573
$ws.Cells.Item(136, 8).Value = 3768.372
$ws.Cells.Item(136, 9).Value = 2205.6072
$ws.Cells.Item(136, 10).Value = 6685.533
$ws.Cells.Item(136, 11).Value = 6616.821599999999
$ws.Cells.Item(136, 12).Value = 20056.599
$ws.Cells.Item(136, 13).Value = -4066.821599999999
$ws.Cells.Item(136, 14).Value = -25156.599
